$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 203.9375
$ws.Range("I15").Value = 203.9375
$ws.Range("K15").Value = 611.8125
$ws.Range("M15").Value = -442.8125

$ws.Range("H33").Value = 637.1177
$ws.Range("I33").Value = 695.2727
$ws.Range("J33").Value = 271.57144
$ws.Range("K33").Value = 695.2727
$ws.Range("L33").Value = 271.57144
$ws.Range("M33").Value = -466.2727
$ws.Range("N33").Value = -729.5714399999999

$ws.Range("H63").Value = 35000
$ws.Range("J63").Value = 35000
$ws.Range("L63").Value = 35000
$ws.Range("N63").Value = -36248

$ws.Range("H66").Value = 35000
$ws.Range("J66").Value = 35000
$ws.Range("L66").Value = 105000
$ws.Range("N66").Value = -111240

$ws.Range("H98").Value = 735.80646
$ws.Range("I98").Value = 671.4815
$ws.Range("J98").Value = 1170
$ws.Range("K98").Value = 671.4815
$ws.Range("L98").Value = 1170
$ws.Range("M98").Value = 826.5185
$ws.Range("N98").Value = -4166

$ws.Range("H122").Value = 735.80646
$ws.Range("I122").Value = 671.4815
$ws.Range("J122").Value = 1170
$ws.Range("K122").Value = 2014.4445
$ws.Range("L122").Value = 3510
$ws.Range("M122").Value = 435.5554999999999
$ws.Range("N122").Value = -8410

$ws.Range("H127").Value = 1074.4166
$ws.Range("I127").Value = 742.2857
$ws.Range("J127").Value = 1539.4
$ws.Range("K127").Value = 2226.8571
$ws.Range("L127").Value = 4618.200000000001
$ws.Range("M127").Value = 2733.1429
$ws.Range("N127").Value = -14538.2

$ws.Range("H129").Value = 6737.4053
$ws.Range("J129").Value = 8486.482
$ws.Range("L129").Value = 25459.446
$ws.Range("N129").Value = -35459.446

$ws.Range("H138").Value = 3127.7188
$ws.Range("I138").Value = 1904.4667
$ws.Range("J138").Value = 3502.1836
$ws.Range("K138").Value = 5713.4001
$ws.Range("L138").Value = 10506.5508
$ws.Range("M138").Value = -573.4000999999998
$ws.Range("N138").Value = -20786.5508

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1007
$ws.Range("I2").Value = 910.9545000000001
$ws.Range("K2").Value = 910.9545000000001
$ws.Range("M2").Value = -797.9545000000001

$ws.Range("H45").Value = 2318.3076
$ws.Range("I45").Value = 2737.3333
$ws.Range("J45").Value = 1959.1428
$ws.Range("K45").Value = 2737.3333
$ws.Range("L45").Value = 1959.1428
$ws.Range("M45").Value = -2360.3333
$ws.Range("N45").Value = -2713.1428

$ws.Range("H116").Value = 1007
$ws.Range("I116").Value = 910.9545000000001
$ws.Range("K116").Value = 910.9545000000001
$ws.Range("M116").Value = 1383.0455

$ws.Range("H122").Value = 1952.2273
$ws.Range("I122").Value = 1638.1082
$ws.Range("K122").Value = 4914.3246
$ws.Range("M122").Value = -2464.3246

$ws.Range("H135").Value = 2517271.5
$ws.Range("J135").Value = 2517271.5
$ws.Range("L135").Value = 2517271.5
$ws.Range("N135").Value = -2527411.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1007
$ws.Range("I3").Value = 910.9545000000001
$ws.Range("K3").Value = 910.9545000000001
$ws.Range("M3").Value = -796.9545000000001

$ws.Range("H99").Value = 1584.1154
$ws.Range("I99").Value = 1411.8667
$ws.Range("J99").Value = 1819
$ws.Range("K99").Value = 1411.8667
$ws.Range("L99").Value = 1819
$ws.Range("M99").Value = 86.13329999999996
$ws.Range("N99").Value = -4815

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2473.566
$ws.Range("I31").Value = 1373.05
$ws.Range("J31").Value = 3140.5454
$ws.Range("K31").Value = 1373.05
$ws.Range("L31").Value = 3140.5454
$ws.Range("M31").Value = -1078.05
$ws.Range("N31").Value = -3730.5454

$ws.Range("H34").Value = 2473.566
$ws.Range("I34").Value = 1373.05
$ws.Range("J34").Value = 3140.5454
$ws.Range("K34").Value = 1373.05
$ws.Range("L34").Value = 3140.5454
$ws.Range("M34").Value = -1171.05
$ws.Range("N34").Value = -3544.5454

$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()

$ws.Range("H122").Value = 7408706
$ws.Range("I122").Value = 13334289
$ws.Range("J122").Value = 1727.5
$ws.Range("K122").Value = 40002867
$ws.Range("L122").Value = 5182.5
$ws.Range("M122").Value = -40000417
$ws.Range("N122").Value = -10082.5

$ws.Range("H132").Value = 1639.0385
$ws.Range("I132").Value = 1366.5238
$ws.Range("J132").Value = 2783.6
$ws.Range("K132").Value = 4099.5714
$ws.Range("L132").Value = 8350.799999999999
$ws.Range("M132").Value = -1569.5714
$ws.Range("N132").Value = -13410.8

$ws.Range("H134").Value = 3380.158
$ws.Range("I134").Value = 3779.4075
$ws.Range("J134").Value = 2400.182
$ws.Range("K134").Value = 11338.2225
$ws.Range("L134").Value = 7200.545999999999
$ws.Range("M134").Value = -8803.2225
$ws.Range("N134").Value = -12270.546

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 4911.8
$ws.Range("I137").Value = 775
$ws.Range("J137").Value = 6416.091
$ws.Range("K137").Value = 2325
$ws.Range("L137").Value = 19248.273
$ws.Range("M137").Value = 2775
$ws.Range("N137").Value = -29448.273

$ws.Range("H140").Value = 1331.6451
$ws.Range("J140").Value = 1522.2
$ws.Range("L140").Value = 4566.6
$ws.Range("N140").Value = -14926.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1565.6757
$ws.Range("I7").Value = 1235.3889
$ws.Range("K7").Value = 1235.3889
$ws.Range("M7").Value = -1123.3889

$ws.Range("H61").Value = 2408.1904
$ws.Range("I61").Value = 2156.9412
$ws.Range("K61").Value = 2156.9412
$ws.Range("M61").Value = -1954.9412

$ws.Range("H113").Value = 2408.1904
$ws.Range("I113").Value = 2156.9412
$ws.Range("K113").Value = 2156.9412
$ws.Range("M113").Value = 13.05879999999979

$ws.Range("H126").Value = 1565.6757
$ws.Range("I126").Value = 1235.3889
$ws.Range("K126").Value = 3706.1667
$ws.Range("M126").Value = -1236.1667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1300.4546
$ws.Range("I122").Value = 1145.1666
$ws.Range("J122").Value = 1999.25
$ws.Range("K122").Value = 3435.4998
$ws.Range("L122").Value = 5997.75
$ws.Range("M122").Value = -985.4998000000001
$ws.Range("N122").Value = -10897.75

$ws.Range("H136").Value = 3494.4285
$ws.Range("I136").Value = 622.129
$ws.Range("J136").Value = 11589.091
$ws.Range("K136").Value = 1866.387
$ws.Range("L136").Value = 34767.273
$ws.Range("M136").Value = 683.6129999999998
$ws.Range("N136").Value = -39867.273
